$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the existing "Late" column (column N),
# pushing Late / heading / Outstanding one column to the right (N->O, O->P, P->Q).
# Match the column width Excel copies from the column immediately to the left (M).
$leftWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab, with K12 selected,
# matching the saved view state in the workbook.
$ws.Activate() | Out-Null
$ws.Range("K12").Select() | Out-Null
